# Updated cryptos list on Sun Sep  3 23:20:38 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) values for the
# crypto list pulled from coinranking.com, and reflects ShibaInu / Litecoin
# swapping ranking places (rows 15 and 16).
#
# Price cells are text (e.g. "214.50", "1.637.02", "0.0\u20857681" — note the
# thousands separators are literal dots and some entries use subscript-digit
# notation for tiny prices), so each Price cell is written with a Text
# NumberFormat first (to stop Excel's automatic "looks like a number" /
# "looks like a date" coercion from rewriting or truncating the string) and
# then has its format cleared again so the cell ends up back on the sheet's
# default (unstyled) look, matching the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript digits used by very small "Price" values (e.g. 0.0\u20857685).
$sub5 = [string][char]8325   # SUBSCRIPT FIVE
$sub8 = [string][char]8328   # SUBSCRIPT EIGHT

# --- Rows 15 & 16 swapped rank: ShibaInu moved up to rank 13 (row 15),
#     Litecoin moved down to rank 14 (row 16), each with refreshed figures.
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0" + $sub5 + "7685"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.56%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.91"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.71%  "

# --- Refresh Price / Volume(1h) figures for the remaining rows.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.952.57"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.53"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("E4").Value = "  -0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5082"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.49%  "

$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2562"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06341"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07755"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.637.22"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5420"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.979.76"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.0000"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.68"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("E20").Value = "  -0.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.887"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.042"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.863"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.20"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1200"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.818"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.234"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04896"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.250"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.167"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.527"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.365"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9105"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.580"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.129.65"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5447"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01561"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.524"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0" + $sub8 + "125"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8092"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.91"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.412"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.775.12"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4524"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.88"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05128"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.34%  "

$ws.Range("E51").Value = "  -0.46%  "
